# Cover page update:
#  - remove one redundant 52-half-pt blank paragraph before the 36-half-pt ones
#  - clear the two "Egyetemi logo kell ide" placeholder runs (paragraph stays, empty)
#  - append " Kovacs Laszlo" to the "Tordeles, grafikai kivitelezes:" line and drop its red color
#  - expand "ISBN" to the full "ISBN 978-963-503-974-6" (two runs)
#  - retarget the two blank paragraphs that followed the second placeholder to match
#    the caps/36-half-pt style used elsewhere on the cover

$d = $word.ActiveDocument

# --- paragraph indices in the ORIGINAL (before-edit) document ---
# 33: blank, caps, sz=52   -> delete (done last so indices below stay valid)
# 36: "Egyetemi logo kell ide" (placeholder #1) -> strip runs, keep empty paragraph
# 60: "Tordeles, grafikai kivitelezes:"          -> drop color, add " Kovacs Laszlo" run
# 63: "ISBN"                                     -> "ISBN 978-963-503-974-" + "6"
# 69: "Egyetemi logo kell ide" (placeholder #2) -> strip runs, keep empty paragraph
# 70: blank, rFonts theme, sz=32                 -> caps, sz=36
# 71: blank, sz=22                               -> caps, sz=36

$blank36Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:caps/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# 1) First placeholder ("Egyetemi logo kell ide") -> remove its two runs, keep pPr.
$p = $d.Paragraphs(36)
$r = $p.Range
$d.Range($r.Start, $r.End - 1).Delete()

# 2) "Tordeles, grafikai kivitelezes:" paragraph -> drop color, add new run.
$p = $d.Paragraphs(60)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>Tördelés, grafikai kivitelezés:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t xml:space="preserve"> Kovács László</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p.Range.InsertXML($xml)

# 3) "ISBN" paragraph -> "ISBN 978-963-503-974-" + "6".
$p = $d.Paragraphs(63)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>ISBN 978-963-503-974-</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="22"/>
      <w:szCs w:val="22"/>
    </w:rPr>
    <w:t>6</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$p.Range.InsertXML($xml)

# 4) Second placeholder ("Egyetemi logo kell ide") -> remove its two runs, keep pPr.
$p = $d.Paragraphs(69)
$r = $p.Range
$d.Range($r.Start, $r.End - 1).Delete()

# 5) The two blank paragraphs right after it become caps/sz=36 (matching the other
#    blank separators used throughout the cover), replacing their old formatting
#    (rFonts theme + sz=32, and sz=22 respectively).
$d.Paragraphs(70).Range.InsertXML($blank36Xml)
$d.Paragraphs(71).Range.InsertXML($blank36Xml)

# 6) Remove one redundant blank 52-half-pt paragraph (of the run of 7) that sits
#    immediately before the two 36-half-pt ones. Done last so it doesn't shift the
#    paragraph numbers used above.
$p = $d.Paragraphs(33)
$p.Range.Delete()

Write-Output "OK"
